$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Insert a new row at position 9 ("filtered_by_ED"), which pushes
#    the old rows 9, 10, 11 down to 10, 11, 12. Excel automatically
#    re-writes the formulas in (what becomes) row 12 to reference the
#    newly-shifted rows (B10-B11 etc.) and updates the dimension.
# ---------------------------------------------------------------
$ws.Rows(9).Insert()

# ---------------------------------------------------------------
# 2. Populate the new row 9 - "filtered_by_ED"
# ---------------------------------------------------------------
$ws.Range("A9").Value = "filtered_by_ED"
$ws.Range("B9").Value = 6079
$ws.Range("C9").Value = 4673
$ws.Range("D9").Value = 12540
$ws.Range("E9").Value = 7207
$ws.Range("F9").Value = 4413
$ws.Range("G9").Value = 2812
$ws.Range("H9").Value = 7302
$ws.Range("I9").Value = 2153

# ---------------------------------------------------------------
# 3. Rename the row that used to be "filtered_by_BR" (now row 10)
# ---------------------------------------------------------------
$ws.Range("A10").Value = "filtered_by_BR(inflection)"

# ---------------------------------------------------------------
# 4. Update the changed values in row 11 (BR_yes_ED_yes)
# ---------------------------------------------------------------
$ws.Range("F11").Value = 2484
$ws.Range("G11").Value = 2037
$ws.Range("K11").Value = 768
$ws.Range("L11").Value = 397

# ---------------------------------------------------------------
# 5. Update K/L values in row 12 (BR_yes_ED_no); B12:I12 formulas were
#    already auto-adjusted by the row insert in step 1.
# ---------------------------------------------------------------
$ws.Range("K12").Value = 120
$ws.Range("L12").Value = 270

# ---------------------------------------------------------------
# 6. Add the new "Summ" column (J) : header + row totals
# ---------------------------------------------------------------
$ws.Range("J1").Value = "Summ"
$ws.Range("J11").Formula = "=SUM(B11:I11)"
$ws.Range("J12").Formula = "=SUM(B12:I12)"
$ws.Range("J12").Style = $ws.Range("I12").Style
$ws.Columns("J").ColumnWidth = 10.5

# ---------------------------------------------------------------
# 7. Row heights for the two rows that got a custom height in Excel
# ---------------------------------------------------------------
$ws.Rows(8).RowHeight = 15
$ws.Rows(9).RowHeight = 15

# ---------------------------------------------------------------
# 8. View: zoom + selection
# ---------------------------------------------------------------
$excel.ActiveWindow.Zoom = 160
$ws.Range("J1").Select()
